$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.348.34'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = '3.509.69'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  +0.05%  '

$c = $ws.Range("D5")
$c.Value = "'591.65"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '

$c = $ws.Range("D6")
$c.Value = "'134.35"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.02%  '

$c = $ws.Range("D9")
$c.Value = "'7.62"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +6.45%  '

$ws.Range("E10").Value = '  -0.07%  '

$ws.Range("E11").Value = '  +3.64%  '

$ws.Range("D12").Value = '4.109.17'
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("E13").Value = '  +1.25%  '

$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("D15").Value = '3.509.76'
$ws.Range("E15").Value = '  +0.27%  '

$c = $ws.Range("D16")
$c.Value = "'25.78"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.95%  '

$ws.Range("D17").Value = '64.337.57'
$ws.Range("E17").Value = '  +0.04%  '

$c = $ws.Range("D18")
$c.Value = "'9.99"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.19%  '

$ws.Range("E19").Value = '  +3.42%  '

$c = $ws.Range("D20")
$c.Value = "'13.55"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.06%  '

$c = $ws.Range("D21")
$c.Value = "'393.87"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.71%  '

$ws.Range("E22").Value = '  +1.08%  '

$ws.Range("D23").Value = '3.651.05'
$ws.Range("E23").Value = '  +0.31%  '

$ws.Range("E24").Value = '  +0.90%  '

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("E26").Value = '  +0.16%  '

$c = $ws.Range("D27")
$c.Value = "'0.0000117"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.75%  '

$ws.Range("E28").Value = '  +0.00%  '

$c = $ws.Range("D29")
$c.Value = "'7.39"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.71%  '

$ws.Range("E30").Value = '  +1.64%  '

$c = $ws.Range("D32")
$c.Value = "'1.46"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -6.46%  '

$c = $ws.Range("D33")
$c.Value = "'0.156"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +6.45%  '

$ws.Range("D34").Value = '3.540.22'
$ws.Range("E34").Value = '  +0.54%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("E36").Value = '  -0.81%  '

$ws.Range("E37").Value = '  +0.65%  '

$c = $ws.Range("D38")
$c.Value = "'6.95"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.58%  '

$ws.Range("E39").Value = '  +0.12%  '

$c = $ws.Range("D40")
$c.Value = "'167.24"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.04%  '

$ws.Range("E41").Value = '  +0.44%  '

$c = $ws.Range("D42")
$c.Value = "'0.812"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.53%  '

$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D44")
$c.Value = "'4.44"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.85%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D45")
$c.Value = "'24.93"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.50%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D46")
$c.Value = "'1.65"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.44%  '

$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range("D47")
$c.Value = "'1.17"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.27%  '

$ws.Range("E48").Value = '  +0.70%  '

$ws.Range("D49").Value = '2.387.22'
$ws.Range("E49").Value = '  -3.33%  '

$c = $ws.Range("D50")
$c.Value = "'0.899"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.81%  '

$ws.Range("E51").Value = '  +0.26%  '
